# Auto-generated edit script: updates Leve profit-calc sheets with refreshed
# market-price data pulled by the scheduled runner (per commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 529.4737
$ws.Range("I2").Value = 528.8823
$ws.Range("K2").Value = 528.8823
$ws.Range("M2").Value = -415.8823
$ws.Range("H4").Value = 241.33333
$ws.Range("I4").Value = 189.6
$ws.Range("K4").Value = 189.6
$ws.Range("M4").Value = -75.59999999999999
$ws.Range("H9").Value = 171.05556
$ws.Range("I9").Value = 171.41667
$ws.Range("K9").Value = 171.41667
$ws.Range("M9").Value = -2.416670000000011
$ws.Range("H133").Value = 77823.71000000001
$ws.Range("J133").Value = 77823.71000000001
$ws.Range("L133").Value = 77823.71000000001
$ws.Range("N133").Value = -87943.71000000001
$ws.Range("H135").Value = 15666.667
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 3500
$ws.Range("K135").Value = 360000
$ws.Range("L135").Value = 31500
$ws.Range("M135").Value = -357465
$ws.Range("N135").Value = -36570
$ws.Range("H137").Value = 7086.026
$ws.Range("I137").Value = 12273.842
$ws.Range("J137").Value = 2157.6
$ws.Range("K137").Value = 36821.526
$ws.Range("L137").Value = 6472.799999999999
$ws.Range("M137").Value = -34271.526
$ws.Range("N137").Value = -11572.8
$ws.Range("H138").Value = 3929.5925
$ws.Range("I138").Value = 1353.6923
$ws.Range("J138").Value = 4746.3413
$ws.Range("K138").Value = 4061.0769
$ws.Range("L138").Value = 14239.0239
$ws.Range("M138").Value = 1078.9231
$ws.Range("N138").Value = -24519.0239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 998
$ws.Range("I4").Value = 998
$ws.Range("K4").Value = 998
$ws.Range("M4").Value = -882
$ws.Range("H5").Value = 835.125
$ws.Range("I5").Value = 476.4
$ws.Range("J5").Value = 1433
$ws.Range("K5").Value = 476.4
$ws.Range("L5").Value = 1433
$ws.Range("M5").Value = -364.4
$ws.Range("N5").Value = -1657
$ws.Range("H50").Value = 3276.5
$ws.Range("J50").Value = 6053
$ws.Range("L50").Value = 6053
$ws.Range("N50").Value = -7481
$ws.Range("H56").Value = 34996.668
$ws.Range("J56").Value = 27495
$ws.Range("L56").Value = 27495
$ws.Range("N56").Value = -28979
$ws.Range("H61").Value = 6581.727
$ws.Range("I61").Value = 8537.375
$ws.Range("K61").Value = 8537.375
$ws.Range("M61").Value = -8325.375
$ws.Range("H132").Value = 5277
$ws.Range("I132").Value = 3924.4
$ws.Range("J132").Value = 6506.636
$ws.Range("K132").Value = 11773.2
$ws.Range("L132").Value = 19519.908
$ws.Range("M132").Value = -9243.200000000001
$ws.Range("N132").Value = -24579.908
$ws.Range("H136").Value = 6581.727
$ws.Range("I136").Value = 8537.375
$ws.Range("K136").Value = 25612.125
$ws.Range("M136").Value = -23062.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 835.125
$ws.Range("I4").Value = 476.4
$ws.Range("J4").Value = 1433
$ws.Range("K4").Value = 476.4
$ws.Range("L4").Value = 1433
$ws.Range("M4").Value = -361.4
$ws.Range("N4").Value = -1663
$ws.Range("H134").Value = 3976.3684
$ws.Range("I134").Value = 3045.4614
$ws.Range("K134").Value = 9136.3842
$ws.Range("M134").Value = -6601.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4346.357
$ws.Range("I31").Value = 1078.4286
$ws.Range("J31").Value = 7614.2856
$ws.Range("K31").Value = 1078.4286
$ws.Range("L31").Value = 7614.2856
$ws.Range("M31").Value = -783.4286
$ws.Range("N31").Value = -8204.285599999999
$ws.Range("H34").Value = 4346.357
$ws.Range("I34").Value = 1078.4286
$ws.Range("J34").Value = 7614.2856
$ws.Range("K34").Value = 1078.4286
$ws.Range("L34").Value = 7614.2856
$ws.Range("M34").Value = -876.4286
$ws.Range("N34").Value = -8018.2856
$ws.Range("H58").Value = 4332.7617
$ws.Range("I58").Value = 3440
$ws.Range("J58").Value = 4689.8667
$ws.Range("K58").Value = 3440
$ws.Range("L58").Value = 4689.8667
$ws.Range("M58").Value = -3237
$ws.Range("N58").Value = -5095.8667
$ws.Range("H99").Value = 247918.28
$ws.Range("I99").Value = 515049.7
$ws.Range("K99").Value = 515049.7
$ws.Range("M99").Value = -513551.7
$ws.Range("H107").Value = 38472280
$ws.Range("I107").Value = 50012904
$ws.Range("K107").Value = 50012904
$ws.Range("M107").Value = -50010984
$ws.Range("H126").Value = 247918.28
$ws.Range("I126").Value = 515049.7
$ws.Range("K126").Value = 1545149.1
$ws.Range("M126").Value = -1542679.1
$ws.Range("H134").Value = 3872.8647
$ws.Range("I134").Value = 3187.4243
$ws.Range("K134").Value = 9562.2729
$ws.Range("M134").Value = -7027.2729
$ws.Range("H136").Value = 4332.7617
$ws.Range("I136").Value = 3440
$ws.Range("J136").Value = 4689.8667
$ws.Range("K136").Value = 10320
$ws.Range("L136").Value = 14069.6001
$ws.Range("M136").Value = -7770
$ws.Range("N136").Value = -19169.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49724344
$ws.Range("I4").Value = 38250836
$ws.Range("K4").Value = 114752508
$ws.Range("M4").Value = -114752396
$ws.Range("H14").Value = 397
$ws.Range("I14").Value = 397
$ws.Range("K14").Value = 1191
$ws.Range("M14").Value = -1018
$ws.Range("H51").Value = 1833
$ws.Range("I51").Value = 862.2727
$ws.Range("K51").Value = 2586.8181
$ws.Range("M51").Value = -2126.8181
$ws.Range("H80").Value = 204624.5
$ws.Range("I80").Value = 9999.333000000001
$ws.Range("J80").Value = 321399.6
$ws.Range("K80").Value = 29997.999
$ws.Range("L80").Value = 964198.7999999999
$ws.Range("M80").Value = -29061.999
$ws.Range("N80").Value = -966070.7999999999
$ws.Range("H83").Value = 204624.5
$ws.Range("I83").Value = 9999.333000000001
$ws.Range("J83").Value = 321399.6
$ws.Range("K83").Value = 89993.997
$ws.Range("L83").Value = 2892596.4
$ws.Range("M83").Value = -85313.997
$ws.Range("N83").Value = -2901956.4
$ws.Range("H113").Value = 3333.3333
$ws.Range("J113").Value = 4249.5
$ws.Range("L113").Value = 12748.5
$ws.Range("N113").Value = -17088.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7199
$ws.Range("I132").Value = 7748.75
$ws.Range("K132").Value = 23246.25
$ws.Range("M132").Value = -20716.25
$ws.Range("H134").Value = 46046.145
$ws.Range("J134").Value = 46046.145
$ws.Range("L134").Value = 138138.435
$ws.Range("N134").Value = -143208.435

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2510.65
$ws.Range("I22").Value = 3132.5454
$ws.Range("J22").Value = 1750.5555
$ws.Range("K22").Value = 3132.5454
$ws.Range("L22").Value = 1750.5555
$ws.Range("M22").Value = -2837.5454
$ws.Range("N22").Value = -2340.5555
$ws.Range("H27").Value = 2510.65
$ws.Range("I27").Value = 3132.5454
$ws.Range("J27").Value = 1750.5555
$ws.Range("K27").Value = 3132.5454
$ws.Range("L27").Value = 1750.5555
$ws.Range("M27").Value = -3025.5454
$ws.Range("N27").Value = -1964.5555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 154000
$ws.Range("I132").Value = 154000
$ws.Range("K132").Value = 462000
$ws.Range("M132").Value = -459470
$ws.Range("H136").Value = 4133.2
$ws.Range("I136").Value = 3916.6667
$ws.Range("J136").Value = 4226
$ws.Range("K136").Value = 11750.0001
$ws.Range("L136").Value = 12678
$ws.Range("M136").Value = -9200.000100000001
$ws.Range("N136").Value = -17778
